$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 258.41349152390097
$ws.Range("C2").Value = 303.32534830605175
$ws.Range("D2").Value = 254.60381929862123
$ws.Range("E2").Value = 305.6639284405972

$ws.Range("B3").Value = 250.20360560472221
$ws.Range("C3").Value = 301.61579614208551
$ws.Range("D3").Value = 248.23193479038781
$ws.Range("E3").Value = 307.6631062766831

$ws.Range("B1:E3").Select() | Out-Null
